$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 10:52"

# Row 8
$ws.Range("B8").Value = 150666
$ws.Range("C8").Value = 18
$ws.Range("E8").Value = 42051

# Row 24
$ws.Range("B24").Value = 15002
$ws.Range("C24").Value = 77
$ws.Range("E24").Value = 2786
$ws.Range("G24").Value = 12
$ws.Range("H24").Value = 522

# Row 51 (country -> Banglades)
$ws.Range("A51").Value = "Banglades"
$ws.Range("B51").Value = 4186
$ws.Range("C51").Value = 414
$ws.Range("D51").Value = 108
$ws.Range("E51").Value = 3951
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 7
$ws.Range("H51").Value = 127

# Row 52 (country -> Finlandia)
$ws.Range("A52").Value = "Finlandia"
$ws.Range("B52").Value = 4129
$ws.Range("D52").Value = 2000
$ws.Range("E52").Value = 1980
$ws.Range("F52").Value = 63
$ws.Range("H52").Value = 149

# Row 71
$ws.Range("B71").Value = 1592
$ws.Range("C71").Value = 33
$ws.Range("D71").Value = 192
$ws.Range("E71").Value = 1355
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 45

# Row 100
$ws.Range("B100").Value = 631
$ws.Range("C100").Value = 19
$ws.Range("D100").Value = 302
$ws.Range("E100").Value = 321
$ws.Range("F100").Value = 11
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 8

# Row 105
$ws.Range("B105").Value = 480
$ws.Range("C105").Value = 6
$ws.Range("D105").Value = 92
$ws.Range("E105").Value = 384

# Row 114
$ws.Range("B114").Value = 334
$ws.Range("C114").Value = 4
$ws.Range("E114").Value = 222

# Row 173 (country -> Malaui)
$ws.Range("A173").Value = "Malaui"
$ws.Range("C173").Value = 10
$ws.Range("D173").Value = 3
$ws.Range("E173").Value = 27
$ws.Range("F173").Value = 1
$ws.Range("H173").Value = 3

# Row 174 (country -> Republica del Chad)
$ws.Range("A174").Value = "Republica del Chad"
$ws.Range("B174").Value = 33
$ws.Range("D174").Value = 8
$ws.Range("E174").Value = 25
$ws.Range("H174").Value = 0

# Row 175 (country -> Guam)
$ws.Range("A175").Value = "Guam"
$ws.Range("B175").Value = 32
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 31

# Row 176 (country -> Suazilandia)
$ws.Range("A176").Value = "Suazilandia"
$ws.Range("B176").Value = 31
$ws.Range("D176").Value = 8
$ws.Range("H176").Value = 1

# Row 177 (country -> Zimbabue)
$ws.Range("A177").Value = "Zimbabue"
$ws.Range("B177").Value = 28
$ws.Range("D177").Value = 2
$ws.Range("E177").Value = 22
$ws.Range("H177").Value = 4

# Row 178 (country -> Angola)
$ws.Range("A178").Value = "Angola"
$ws.Range("B178").Value = 25
$ws.Range("D178").Value = 6
$ws.Range("E178").Value = 17
$ws.Range("F178").Value = 0
$ws.Range("H178").Value = 2

# Row 179 (country -> Antigua y Barbuda)
$ws.Range("A179").Value = "Antigua y Barbuda"
$ws.Range("B179").Value = 24
$ws.Range("D179").Value = 10
$ws.Range("E179").Value = 11
$ws.Range("F179").Value = 1
$ws.Range("H179").Value = 3

# Row 180 (country -> Timor Oriental)
$ws.Range("A180").Value = "Timor Oriental"
$ws.Range("D180").Value = 1
$ws.Range("E180").Value = 22
$ws.Range("F180").Value = 0
$ws.Range("H180").Value = 0
